$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,C,D,E,F,I,M,N,O across rows 2-25 (data rows 0-23)
$colIndices = @(2, 3, 4, 5, 6, 9, 13, 14, 15)

$data = @(
    @(0.6194284934352936, 0.1725593584462786, 0.0135231891998373, 0.1282636622105215, 0.5034795799655143, 0.3408651183409077, 0.3083193749590762, 0.9311592783408713, 1.612172864956648),
    @(0.5420733954248078, 0.1519890321154094, 0.01213268093606956, 0.1218239920399995, 0.4948761282577863, 0.3436564154520951, 0.2744500294856138, 0.9400192022935059, 1.596948262432562),
    @(0.4944770899284379, 0.1392938856775459, 0.01127379038392462, 0.1179836283755336, 0.4899921034708754, 0.3456378482473923, 0.2537132601275118, 0.9458704867437859, 1.58889822309817),
    @(0.4750571129632419, 0.134104397221904, 0.01092252354636258, 0.1164469638192074, 0.4881018492000919, 0.34651243347205, 0.2452777522453076, 0.948358372953173, 1.5859434156595),
    @(0.4718310159394719, 0.1332417193931406, 0.01086412043557061, 0.1161935045138378, 0.4877940097165805, 0.3466617082070336, 0.2438779428478526, 0.9487777329735181, 1.58547241930836),
    @(0.4942152815112308, 0.1392239633949544, 0.01126905815480939, 0.1179627900721911, 0.4899662060916299, 0.3456493715760871, 0.2535994355825011, 0.945903620384918, 1.588857055957931),
    @(0.5927778029779063, 0.1654802915242328, 0.01304481692855575, 0.1260195259336783, 0.5004302778597634, 0.3417719362566984, 0.2966289509915043, 0.9341288392769584, 1.606653590733941),
    @(0.7852290036541376, 0.2164475552774263, 0.01648562158362665, 0.1427327368123628, 0.5241227882272028, 0.336298087261941, 0.3814824730977904, 0.9143006215767926, 1.651888208096096),
    @(0.9260818240026651, 0.2535701757806237, 0.01898735174597732, 0.1555876722083909, 0.5434807257051233, 0.3335844847410456, 0.4441237877434077, 0.901720308371992, 1.691480876081044),
    @(0.9900352928326583, 0.2703871859748972, 0.02011956918149593, 0.1615648086037638, 0.552714613214647, 0.3326361203143371, 0.4726887606742025, 0.8964286139695119, 1.710886183006835),
    @(1.014234513232566, 0.2767450911489675, 0.02054745068963371, 0.1638470985907645, 0.5562730213843992, 0.3323183020425979, 0.4835155948606769, 0.8944867908889478, 1.718435883464736),
    @(1.009023622240875, 0.275376265651488, 0.02045533762079543, 0.1633547225914853, 0.5555039055048923, 0.332384910049587, 0.4811834027971287, 0.8949022380919729, 1.716800953093468),
    @(0.9920265559777022, 0.2709104630305319, 0.02015478872873189, 0.1617521940440696, 0.5530061271349922, 0.3326091443882788, 0.4735792926846045, 0.8962676157349279, 1.711503261737676),
    @(0.9816129098119291, 0.2681736777143726, 0.01997058047859923, 0.160773065736656, 0.5514842117100898, 0.3327518788152659, 0.4689228435886434, 0.8971120275688591, 1.708284516414409),
    @(0.9218997779911433, 0.2524697113040304, 0.01891323909109843, 0.1551996789619281, 0.5428858970606427, 0.3336522335804801, 0.4422583870204733, 0.9020748111944599, 1.690240822609553),
    @(0.8852358523258772, 0.2428176792677448, 0.01826308212041283, 0.1518139155644036, 0.5377208449864668, 0.3342779621945731, 0.4259182816420264, 0.9052297649853216, 1.67952932413553),
    @(0.8641363936358744, 0.2372594865412339, 0.01788858162592533, 0.1498786736076809, 0.5347903038725974, 0.334664780529323, 0.4165263939940189, 0.9070849912034546, 1.673499552609684),
    @(0.8569905669909303, 0.2353764515522414, 0.01776168897001185, 0.1492255139549528, 0.5338049813531143, 0.334800368180435, 0.4133475813613501, 0.9077201078600794, 1.671480484180449),
    @(0.889139971578004, 0.2438458392075802, 0.0183323492823817, 0.152173075388724, 0.5382665052964768, 0.3342085653944551, 0.4276570404212663, 0.904889714678454, 1.680655996881853),
    @(0.9970195194483722, 0.2722224597310969, 0.02024309089736676, 0.1622223809726151, 0.5537381080496999, 0.3325421589150004, 0.4758125358242182, 0.8958648875950743, 1.713053851507112),
    @(1.067416269712737, 0.2907078952989366, 0.02148681906179917, 0.1689003481057227, 0.5642096319423757, 0.3316939123615477, 0.5073426153060865, 0.8903281966029439, 1.735401616727017),
    @(1.029854528956491, 0.2808474657763895, 0.02082348876265883, 0.1653260198407551, 0.5585877795687111, 0.3321245437193348, 0.4905091542973281, 0.893250140244632, 1.723366512656298),
    @(0.8873749842153416, 0.2433810364620399, 0.01830103582818765, 0.1520106641235799, 0.5380196910645481, 0.3342398553511963, 0.4268709405496764, 0.9050433224582477, 1.680146228196776),
    @(0.7332583334284095, 0.2027159123626916, 0.01555932986391895, 0.13811160773583, 0.5173719367047838, 0.3375499010877867, 0.3584755759702745, 0.9193155903011885, 1.638538550520138)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $colIndices.Count; $j++) {
        $ws.Cells.Item($rowNum, $colIndices[$j]).Value = $rowValues[$j]
    }
}
